$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns D ("Rage Increase") and E ("Impact") ---
$ws.Range("D1").Value = "Rage Increase"
$ws.Range("E1").Value = "Impact"

# --- Data rows 2-8: Rage Increase (D) / Impact (E) values ---
$rageIncrease = @{2=1; 3=1; 4=1; 5=1; 6=5; 7=2; 8=2}
$impact       = @{2=1; 3=1; 4=1; 5=1; 6=2; 7=1; 8=1}

foreach ($r in 2..8) {
    $ws.Cells.Item($r, 4).Value = $rageIncrease[$r]
    $ws.Cells.Item($r, 5).Value = $impact[$r]
}

# --- Copy the existing header/data formatting (style s="2") onto the new cells ---
$ws.Range("C1:C8").Copy()
$ws.Range("D1:E8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column widths: D gets an explicit best-fit width, matching column B (11.5 chars) ---
$ws.Columns.Item(4).ColumnWidth = 11.5 - 5/7

# --- Selection, matching the saved workbook state ---
$ws.Range("F8").Select()
